# ---------------------------------------------------------------------------
# Applies the commit's changes to PlayerPerformance_3863.xlsx:
#   1. Inserts a brand-new "Player Info" worksheet as the first sheet, with
#      player ID / NAME / BATTING_HAND / BOWL_STYLE data.
#   2. Renames the MATCH_CARD_LINK column on "ODI Batting" to MATCH_CODE and
#      replaces the full scorecard URLs with the bare numeric match code.
#   3. Does the same on "ODI Bowling" (its MATCH_CARD_LINK column is column B).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "Player Info" sheet in front of "ODI Batting"
# ---------------------------------------------------------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")

$infoSheet = $wb.Worksheets.Add($battingSheet)
$infoSheet.Name = "Player Info"

$headers = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($c = 1; $c -le $headers.Count; $c++) {
    $infoSheet.Cells.Item(1, $c).Value = $headers[$c - 1]
}

$infoSheet.Cells.Item(2, 1).NumberFormat = "@"
$infoSheet.Cells.Item(2, 1).Value = "3863"
$infoSheet.Cells.Item(2, 2).Value = "Kraigg Clairmonte Brathwaite"
$infoSheet.Cells.Item(2, 3).Value = "Right Handed"
$infoSheet.Cells.Item(2, 4).Value = "Right Arm Off Break"

# header formatting to match the look of the other sheets' header rows
$headerRange = $infoSheet.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

$infoSheet.Range("A1").Select()

# ---------------------------------------------------------------------------
# 2. "ODI Batting": MATCH_CARD_LINK (col D) -> MATCH_CODE, URL -> bare code
#    (re-fetch the sheet by name since indices shifted after the insert)
# ---------------------------------------------------------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$battingSheet.Range("D1").Value = "MATCH_CODE"

$lastRow = $battingSheet.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $battingSheet.Cells.Item($r, 4)
    $url = $cell.Text
    $code = $url -replace ".*MatchCode=", ""
    $cell.NumberFormat = "@"
    $cell.Value = $code
}

# ---------------------------------------------------------------------------
# 3. "ODI Bowling": MATCH_CARD_LINK (col B) -> MATCH_CODE, URL -> bare code
# ---------------------------------------------------------------------------
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$bowlingSheet.Range("B1").Value = "MATCH_CODE"

$lastRow = $bowlingSheet.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $bowlingSheet.Cells.Item($r, 2)
    $url = $cell.Text
    $code = $url -replace ".*MatchCode=", ""
    $cell.NumberFormat = "@"
    $cell.Value = $code
}

$wb.Save()
